# Update the "Bill Summary" sheet of the contractor bill with the new
# quantities / rates / item descriptions.
#
# Notes on quoting:
#  - Columns C, F (and H on the data rows) hold true numbers, so those are
#    assigned bare numeric literals.
#  - Columns A, D, E, G (and H on the total rows) are stored as text even
#    when the text looks like a number (e.g. "7.0", "350.00"). Assigning a
#    plain numeric-looking string to a Range.Value lets Excel re-interpret
#    it as a real number and drop the trailing zeros, so those values are
#    prefixed with a leading apostrophe ('), Excel's standard "treat this
#    as text" marker, to preserve the exact text (matches how a user
#    typing '350.00 into a cell keeps it as literal text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("C8").Value = 7
$ws.Range("G8").Value = '''350.00'

# Row 9
$ws.Range("A9").Value = 'Each'
$ws.Range("C9").Value = 39
$ws.Range("D9").Value = '''7.0'
$ws.Range("E9").Value = 'Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F9").Value = 30
$ws.Range("G9").Value = '''1170.00'

# Row 10
$ws.Range("A10").Value = 'Each'
$ws.Range("C10").Value = 75
$ws.Range("D10").Value = '''8.0'
$ws.Range("E10").Value = 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F10").Value = 30
$ws.Range("G10").Value = '''2250.00'

# Row 11
$ws.Range("A11").Value = 'Each'
$ws.Range("C11").Value = 52
$ws.Range("D11").Value = '''10.0'
$ws.Range("E11").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F11").Value = 303
$ws.Range("G11").Value = '''15756.00'

# Row 12
$ws.Range("C12").Value = 62
$ws.Range("D12").Value = '''36'
$ws.Range("E12").Value = 'Total'

# Row 14 (Grand Total)
$ws.Range("G14").Value = '''19526.00'
$ws.Range("H14").Value = '''19526.00'

# Row 16 (Net payable amount)
$ws.Range("G16").Value = '''19526.00'
$ws.Range("H16").Value = '''19526.00'
